$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 42.75280866666666
$ws.Range("H2").Value = 128.258426
$ws.Range("I2").Value = 0.8529286054750734
$ws.Range("J2").Value = 0.8529286054750735
$ws.Range("M2").Value = 1.306600666666667
$ws.Range("N2").Value = 3.919802
$ws.Range("O2").Value = 0.8137131711319011
$ws.Range("P2").Value = 0.8137131711319011
$ws.Range("Q2").Value = 55.8608483057391
$ws.Range("R2").Value = 502.7476347516519
$ws.Range("S2").Value = 0.6940392403102322
$ws.Range("T2").Value = 0.6940392403102323

# Row 3
$ws.Range("G3").Value = 42.75280866666666
$ws.Range("H3").Value = 128.258426
$ws.Range("I3").Value = 0.8529286054750734
$ws.Range("J3").Value = 0.8529286054750735
$ws.Range("O3").Value = 0.006864390964089149
$ws.Range("P3").Value = 0.006864390964089149
$ws.Range("Q3").Value = 0.4712357080602221
$ws.Range("R3").Value = 4.241121372542
$ws.Range("S3").Value = 0.005854835412436253
$ws.Range("T3").Value = 0.005854835412436254

# Row 4
$ws.Range("G4").Value = 42.75280866666666
$ws.Range("H4").Value = 128.258426
$ws.Range("I4").Value = 0.8529286054750734
$ws.Range("J4").Value = 0.8529286054750735
$ws.Range("M4").Value = 0.2881033333333333
$ws.Range("N4").Value = 0.86431
$ws.Range("O4").Value = 0.1794224379040098
$ws.Range("P4").Value = 0.1794224379040098
$ws.Range("Q4").Value = 12.31722668622889
$ws.Range("R4").Value = 110.85504017606
$ws.Range("S4").Value = 0.153034529752405
$ws.Range("T4").Value = 0.1530345297524051

# Row 5
$ws.Range("I5").Value = 0.04642608686423023
$ws.Range("J5").Value = 0.04642608686423023
$ws.Range("M5").Value = 1.306600666666667
$ws.Range("N5").Value = 3.919802
$ws.Range("O5").Value = 0.8137131711319011
$ws.Range("P5").Value = 0.8137131711319011
$ws.Range("Q5").Value = 3.040583442863111
$ws.Range("R5").Value = 27.365250985768
$ws.Range("S5").Value = 0.03777751836553788
$ws.Range("T5").Value = 0.03777751836553788

# Row 6
$ws.Range("I6").Value = 0.04642608686423023
$ws.Range("J6").Value = 0.04642608686423023
$ws.Range("O6").Value = 0.006864390964089149
$ws.Range("P6").Value = 0.006864390964089149
$ws.Range("S6").Value = 0.0003186868111688399
$ws.Range("T6").Value = 0.0003186868111688399

# Row 7
$ws.Range("I7").Value = 0.04642608686423023
$ws.Range("J7").Value = 0.04642608686423023
$ws.Range("M7").Value = 0.2881033333333333
$ws.Range("N7").Value = 0.86431
$ws.Range("O7").Value = 0.1794224379040098
$ws.Range("P7").Value = 0.1794224379040098
$ws.Range("Q7").Value = 0.6704437304488889
$ws.Range("R7").Value = 6.03399357404
$ws.Range("S7").Value = 0.008329881687523513
$ws.Range("T7").Value = 0.008329881687523515

# Row 8
$ws.Range("G8").Value = 5.044817999999999
$ws.Range("I8").Value = 0.1006453076606963
$ws.Range("J8").Value = 0.1006453076606963
$ws.Range("M8").Value = 1.306600666666667
$ws.Range("N8").Value = 3.919802
$ws.Range("O8").Value = 0.8137131711319011
$ws.Range("P8").Value = 0.8137131711319011
$ws.Range("Q8").Value = 6.591562562011998
$ws.Range("R8").Value = 59.32406305810799
$ws.Range("S8").Value = 0.081896412456131
$ws.Range("T8").Value = 0.081896412456131

# Row 9
$ws.Range("G9").Value = 5.044817999999999
$ws.Range("I9").Value = 0.1006453076606963
$ws.Range("J9").Value = 0.1006453076606963
$ws.Range("O9").Value = 0.006864390964089149
$ws.Range("P9").Value = 0.006864390964089149
$ws.Range("R9").Value = 0.5004509904179999
$ws.Range("S9").Value = 0.0006908687404840561
$ws.Range("T9").Value = 0.0006908687404840561

# Row 10
$ws.Range("G10").Value = 5.044817999999999
$ws.Range("I10").Value = 0.1006453076606963
$ws.Range("J10").Value = 0.1006453076606963
$ws.Range("M10").Value = 0.2881033333333333
$ws.Range("N10").Value = 0.86431
$ws.Range("O10").Value = 0.1794224379040098
$ws.Range("P10").Value = 0.1794224379040098
$ws.Range("S10").Value = 0.01805802646408124
$ws.Range("T10").Value = 0.01805802646408125
